$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text storage first so numeric-looking strings
# (e.g. "0.999", "66.771.42") are kept as literal text like the source file,
# matching the original inlineStr cells instead of being parsed as numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '66.771.42'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '3.213.60'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.53%  '
$ws.Range('D5').Value = '579.32'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').Value = '139.65'
$ws.Range('E6').Value = '  -7.83%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -1.19%  '
$ws.Range('D8').Value = '3.204.26'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -2.74%  '
$ws.Range('D11').Value = '6.27'
$ws.Range('E11').Value = '  +9.84%  '
$ws.Range('D12').Value = '0.474'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = '35.52'
$ws.Range('E14').Value = '  -3.50%  '
$ws.Range('D15').Value = '3.719.66'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '66.683.75'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '3.214.08'
$ws.Range('E17').Value = '  -2.37%  '
$ws.Range('E18').Value = '  -3.18%  '
$ws.Range('D19').Value = '500.51'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '14.13'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').Value = '0.709'
$ws.Range('E22').Value = '  -4.61%  '
$ws.Range('D23').Value = '7.32'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').Value = '81.36'
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('D25').Value = '12.75'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').Value = '3.06'
$ws.Range('E27').Value = '  -4.76%  '
$ws.Range('D28').Value = '2.02'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '27.61'
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  +6.30%  '
$ws.Range('E31').Value = '  +5.95%  '
$ws.Range('D32').Value = '2.49'
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('D34').Value = '495.53'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '53.98'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '6.00'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '5.24'
$ws.Range('E37').Value = '  -3.72%  '
$ws.Range('D38').Value = '0.0407'
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('D39').Value = '0.0805'
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('D40').Value = '8.44'
$ws.Range('E40').Value = '  -5.96%  '
$ws.Range('E41').Value = '  +8.50%  '
$ws.Range('D42').Value = '2.829.29'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').Value = '2.50'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = '0.248'
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('D46').Value = '24.75'
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('D47').Value = '120.21'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('B48').Value = 'PEPE'
$ws.Range('C48').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D48').Value = '0.0₃0524'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '1.99'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').Value = '0.108'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').Value = '2.08'
$ws.Range('E51').Value = '  -12.67%  '

# Restore the default (unstyled) look so we do not leave an explicit
# Text number-format behind on cells that never had one originally.
$dataRange.Style = "Normal"
